# Apply the "corrections in target metrics" edit:
#  - Insert 3 new rows after current row 19 (new CK19 - A - 22/23/24 images)
#  - Insert 3 new rows after current row 37 (new CK19 - B - 22/23/24 images)
#  - Update the selection on "per image"
#  - Update the AVERAGE/AVERAGEIF formulas on "per category" to cover the
#    extended range (B2:B295 / C2:C295 instead of B2:B289 / C2:C289)

$wb = $excel.ActiveWorkbook
$wsImage = $wb.Worksheets.Item("per image")
$wsCategory = $wb.Worksheets.Item("per category")

# --- Insert the first block (lower row numbers) first, so the brand-new
#     shared strings get appended in the same order as the target file
#     (CK19 - A - 22/23/24, then CK19 - B - 22/23/24). ---

# New rows for "CK19 - A - 22/23/24", inserted right before old row 20
# (i.e. right after the row currently holding "CK19 - A - 21").
$wsImage.Rows.Item(20).Insert()
$wsImage.Range("A20").Value = "CK19 - A - 22"
$wsImage.Range("B20").Value = 0.666
$wsImage.Range("C20").Value = "ConA"

$wsImage.Rows.Item(21).Insert()
$wsImage.Range("A21").Value = "CK19 - A - 23"
$wsImage.Range("B21").Value = 0.264
$wsImage.Range("C21").Value = "OVA"

$wsImage.Rows.Item(22).Insert()
$wsImage.Range("A22").Value = "CK19 - A - 24"
$wsImage.Range("B22").Value = 0.193
$wsImage.Range("C22").Value = "Unstimulated"

# New rows for "CK19 - B - 22/23/24", inserted right before what is now
# row 41 (i.e. right after the row currently holding "CK19 - B - 21",
# which was old row 37 and shifted down to row 40 by the inserts above).
$wsImage.Rows.Item(41).Insert()
$wsImage.Range("A41").Value = "CK19 - B - 22"
$wsImage.Range("B41").Value = 0.465
$wsImage.Range("C41").Value = "ConA"

$wsImage.Rows.Item(42).Insert()
$wsImage.Range("A42").Value = "CK19 - B - 23"
$wsImage.Range("B42").Value = 0.297
$wsImage.Range("C42").Value = "OVA"

$wsImage.Rows.Item(43).Insert()
$wsImage.Range("A43").Value = "CK19 - B - 24"
$wsImage.Range("B43").Value = 0.209
$wsImage.Range("C43").Value = "Unstimulated"

# Move the active selection on "per image" to C41 (matches the edit target).
$wsImage.Range("C41").Select()

# Update the dependent formulas on "per category" to use the new row count.
$wsCategory.Range("B1").Formula = "=AVERAGE('per image'!B2:'per image'!B295)"
$wsCategory.Range("B2").Formula = "=AVERAGEIF('per image'!C2:'per image'!C295, ""Unstimulated"", 'per image'!B2:'per image'!B295)"
$wsCategory.Range("B3").Formula = "=AVERAGEIF('per image'!C2:'per image'!C295, ""OVA"", 'per image'!B2:'per image'!B295)"
$wsCategory.Range("B4").Formula = "=AVERAGEIF('per image'!C2:'per image'!C295, ""ConA"", 'per image'!B2:'per image'!B295)"

Write-Host "edit applied"
